$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.41908946361288
$ws.Range("C2").Value = 7.626184542585443
$ws.Range("D2").Value = 14.53761406220659
$ws.Range("E2").Value = 15.753816265174
$ws.Range("G2").Value = 30.06075141182773
$ws.Range("H2").Value = 14.53618689881979
$ws.Range("I2").Value = 20.18342821128782
$ws.Range("J2").Value = 9.23275033487905
$ws.Range("M2").Value = 17.17873311764749
$ws.Range("N2").Value = 17.39541469407961
$ws.Range("O2").Value = 22.31859631746758
$ws.Range("B3").Value = 11.93344286615437
$ws.Range("C3").Value = 7.206989061481051
$ws.Range("D3").Value = 14.5219544028458
$ws.Range("E3").Value = 15.77144270344167
$ws.Range("G3").Value = 30.03187385834743
$ws.Range("H3").Value = 14.57642666397888
$ws.Range("I3").Value = 20.27646498285445
$ws.Range("J3").Value = 9.254521919266084
$ws.Range("M3").Value = 17.03427758976533
$ws.Range("N3").Value = 17.44196402964048
$ws.Range("O3").Value = 22.36616407340034
$ws.Range("B4").Value = 11.62608769536219
$ws.Range("C4").Value = 6.936858536567611
$ws.Range("D4").Value = 14.51534714894636
$ws.Range("E4").Value = 15.78520815651824
$ws.Range("G4").Value = 30.02451225307438
$ws.Range("H4").Value = 14.60373969868792
$ws.Range("I4").Value = 20.33785371014577
$ws.Range("J4").Value = 9.268781138643474
$ws.Range("M4").Value = 16.94763064336586
$ws.Range("N4").Value = 17.47227168037263
$ws.Range("O4").Value = 22.40071933365079
$ws.Range("B5").Value = 11.49871613129772
$ws.Range("C5").Value = 6.823667463309212
$ws.Range("D5").Value = 14.51341362888146
$ws.Range("E5").Value = 15.79155742015599
$ws.Range("G5").Value = 30.02412007244627
$ws.Range("H5").Value = 14.61552450835436
$ws.Range("I5").Value = 20.3639408476037
$ws.Range("J5").Value = 9.274816412403696
$ws.Range("M5").Value = 16.91286838421974
$ws.Range("N5").Value = 17.48505720807279
$ws.Range("O5").Value = 22.41614256108532
$ws.Range("B6").Value = 11.47744378693007
$ws.Range("C6").Value = 6.804687404565072
$ws.Range("D6").Value = 14.51313847740805
$ws.Range("E6").Value = 15.79265637717078
$ws.Range("G6").Value = 30.02421240674892
$ws.Range("H6").Value = 14.61752087794449
$ws.Range("I6").Value = 20.3683372254074
$ws.Range("J6").Value = 9.27583213711241
$ws.Range("M6").Value = 16.90713008882565
$ws.Range("N6").Value = 17.48720652832436
$ws.Range("O6").Value = 22.41878450894438
$ws.Range("B7").Value = 11.62437825190647
$ws.Range("C7").Value = 6.935344459690161
$ws.Range("D7").Value = 14.51531799659251
$ws.Range("E7").Value = 15.78529079024775
$ws.Range("G7").Value = 30.02449640705509
$ws.Range("H7").Value = 14.60389598354868
$ws.Range("I7").Value = 20.33820119692565
$ws.Range("J7").Value = 9.268861622850659
$ws.Range("M7").Value = 16.94715957105893
$ws.Range("N7").Value = 17.47244234838175
$ws.Range("O7").Value = 22.400921908838
$ws.Range("B8").Value = 12.25364442129866
$ws.Range("C8").Value = 7.484341250579252
$ws.Range("D8").Value = 14.53159207353902
$ws.Range("E8").Value = 15.75928300101691
$ws.Range("G8").Value = 30.0486432615013
$ws.Range("H8").Value = 14.54952036453097
$ws.Range("I8").Value = 20.21462148128971
$ws.Range("J8").Value = 9.240072401041157
$ws.Range("M8").Value = 17.12851883805358
$ws.Range("N8").Value = 17.41110715008495
$ws.Range("O8").Value = 22.33388575988436
$ws.Range("B9").Value = 13.40788650208441
$ws.Range("C9").Value = 8.456565776441856
$ws.Range("D9").Value = 14.58723364555084
$ws.Range("E9").Value = 15.73163558889445
$ws.Range("G9").Value = 30.17813495108189
$ws.Range("H9").Value = 14.46359429262575
$ws.Range("I9").Value = 20.0061773841696
$ws.Range("J9").Value = 9.190672338667431
$ws.Range("M9").Value = 17.49898981280507
$ws.Range("N9").Value = 17.30448701803384
$ws.Range("O9").Value = 22.24499690289855
$ws.Range("B10").Value = 14.1992264737897
$ws.Range("C10").Value = 9.1038783248021
$ws.Range("D10").Value = 14.64236284620534
$ws.Range("E10").Value = 15.72555612062308
$ws.Range("G10").Value = 30.32298304345036
$ws.Range("H10").Value = 14.41312518435508
$ws.Range("I10").Value = 19.87378557361394
$ws.Range("J10").Value = 9.158655654736378
$ws.Range("M10").Value = 17.77824723647113
$ws.Range("N10").Value = 17.23442663331137
$ws.Range("O10").Value = 22.20580015304385
$ws.Range("B11").Value = 14.54555993938972
$ws.Range("C11").Value = 9.383329444870519
$ws.Range("D11").Value = 14.67047967691755
$ws.Range("E11").Value = 15.72587567174014
$ws.Range("G11").Value = 30.39953571672728
$ws.Range("H11").Value = 14.39292211026277
$ws.Range("I11").Value = 19.81808299965159
$ws.Range("J11").Value = 9.145014301230082
$ws.Range("M11").Value = 17.90639743424394
$ws.Range("N11").Value = 17.2043399168569
$ws.Range("O11").Value = 22.19366395392165
$ws.Range("B12").Value = 14.67465119521562
$ws.Range("C12").Value = 9.486961939149277
$ws.Range("D12").Value = 14.68155798860531
$ws.Range("E12").Value = 15.72643942291301
$ws.Range("G12").Value = 30.430041215098
$ws.Range("H12").Value = 14.38566851900439
$ws.Range("I12").Value = 19.79764209657048
$ws.Range("J12").Value = 9.139981061557695
$ws.Range("M12").Value = 17.95504872991522
$ws.Range("N12").Value = 17.19320262293185
$ws.Range("O12").Value = 22.1898885455952
$ws.Range("B13").Value = 14.64694196813113
$ws.Range("C13").Value = 9.464740622226184
$ws.Range("D13").Value = 14.67915300184962
$ws.Range("E13").Value = 15.72629833469066
$ws.Range("G13").Value = 30.42340414514748
$ws.Range("H13").Value = 14.3872130492123
$ws.Range("I13").Value = 19.80201536003713
$ws.Range("J13").Value = 9.141059174135592
$ws.Range("M13").Value = 17.94456588883592
$ws.Range("N13").Value = 17.19558986783881
$ws.Range("O13").Value = 22.19066514600735
$ws.Range("B14").Value = 14.55622208958036
$ws.Range("C14").Value = 9.391899385195646
$ws.Range("D14").Value = 14.67138249202459
$ws.Range("E14").Value = 15.72591318462767
$ws.Range("G14").Value = 30.40201512945673
$ws.Range("H14").Value = 14.39231739560945
$ws.Range("I14").Value = 19.81638822636267
$ws.Range("J14").Value = 9.144597560937292
$ws.Range("M14").Value = 17.91039773830259
$ws.Range("N14").Value = 17.20341852021802
$ws.Range("O14").Value = 22.1933369020528
$ws.Range("B15").Value = 14.50038298653215
$ws.Range("C15").Value = 9.346996089572066
$ws.Range("D15").Value = 14.66667878785757
$ws.Range("E15").Value = 15.72573489673542
$ws.Range("G15").Value = 30.38911072291921
$ws.Range("H15").Value = 14.39549565752515
$ws.Range("I15").Value = 19.82527705152048
$ws.Range("J15").Value = 9.146782164144609
$ws.Range("M15").Value = 17.88948376842053
$ws.Range("N15").Value = 17.20824710139993
$ws.Range("O15").Value = 22.19508029267152
$ws.Range("B16").Value = 14.17630966898561
$ws.Range("C16").Value = 9.085310635649837
$ws.Range("D16").Value = 14.6405859972494
$ws.Range("E16").Value = 15.72559724709353
$ws.Range("G16").Value = 30.31819333243736
$ws.Range("H16").Value = 14.41450100680676
$ws.Range("I16").Value = 19.87751705765598
$ws.Range("J16").Value = 9.159565699756257
$ws.Range("M16").Value = 17.76989161111002
$ws.Range("N16").Value = 17.23642871489361
$ws.Range("O16").Value = 22.20670800467453
$ws.Range("B17").Value = 13.97393310147571
$ws.Range("C17").Value = 8.920905490758537
$ws.Range("D17").Value = 14.62535311412691
$ws.Range("E17").Value = 15.72630242106145
$ws.Range("G17").Value = 30.27740750708977
$ws.Range("H17").Value = 14.42686640557629
$ws.Range("I17").Value = 19.91072452570465
$ws.Range("J17").Value = 9.167644215637091
$ws.Range("M17").Value = 17.69678505801502
$ws.Range("N17").Value = 17.25417369820466
$ws.Range("O17").Value = 22.21530091344039
$ws.Range("B18").Value = 13.85625130005825
$ws.Range("C18").Value = 8.824932077785405
$ws.Range("D18").Value = 14.61687792187588
$ws.Range("E18").Value = 15.72699847459336
$ws.Range("G18").Value = 30.25495315472119
$ws.Range("H18").Value = 14.43423798984692
$ws.Range("I18").Value = 19.9302501971129
$ws.Range("J18").Value = 9.172377678127303
$ws.Range("M18").Value = 17.65484320955169
$ws.Range("N18").Value = 17.26454810967969
$ws.Range("O18").Value = 22.22077927253232
$ws.Range("B19").Value = 13.81618950500591
$ws.Range("C19").Value = 8.792195676759848
$ws.Range("D19").Value = 14.61405771878773
$ws.Range("E19").Value = 15.72728405572307
$ws.Range("G19").Value = 30.24752348536976
$ws.Range("H19").Value = 14.43677840249046
$ws.Range("I19").Value = 19.93693429222602
$ws.Range("J19").Value = 9.173995284249637
$ws.Range("M19").Value = 17.64066193466221
$ws.Range("N19").Value = 17.26808957597566
$ws.Range("O19").Value = 22.22272614896361
$ws.Range("B20").Value = 13.99560965215879
$ws.Range("C20").Value = 8.93855307085261
$ws.Range("D20").Value = 14.62694508397485
$ws.Range("E20").Value = 15.72619729867351
$ws.Range("G20").Value = 30.2816453673872
$ws.Range("H20").Value = 14.42552324295074
$ws.Range("I20").Value = 19.90714546994658
$ws.Range("J20").Value = 9.166775250833469
$ws.Range("M20").Value = 17.70455655612036
$ws.Range("N20").Value = 17.25226733537805
$ws.Range("O20").Value = 22.21433070349741
$ws.Range("B21").Value = 14.58292525162825
$ws.Range("C21").Value = 9.413354255674642
$ws.Range("D21").Value = 14.67365322874531
$ws.Range("E21").Value = 15.72601430492816
$ws.Range("G21").Value = 30.40825658367068
$ws.Range("H21").Value = 14.39080734957962
$ws.Range("I21").Value = 19.81214884336809
$ws.Range("J21").Value = 9.143554658890221
$ws.Range("M21").Value = 17.92043068981109
$ws.Range("N21").Value = 17.20111211554605
$ws.Range("O21").Value = 22.19252987049158
$ws.Range("B22").Value = 14.95474117733383
$ws.Range("C22").Value = 9.710892201739044
$ws.Range("D22").Value = 14.70668929313313
$ws.Range("E22").Value = 15.72847477143302
$ws.Range("G22").Value = 30.49983511553274
$ws.Range("H22").Value = 14.37043204364254
$ws.Range("I22").Value = 19.75386701751297
$ws.Range("J22").Value = 9.129150509476254
$ws.Range("M22").Value = 18.06222042328829
$ws.Range("N22").Value = 17.16917046872756
$ws.Range("O22").Value = 22.18306348340292
$ws.Range("B23").Value = 14.75742457445408
$ws.Range("C23").Value = 9.5532675933171
$ws.Range("D23").Value = 14.68882975211435
$ws.Range("E23").Value = 15.72692585686603
$ws.Range("G23").Value = 30.45015600657556
$ws.Range("H23").Value = 14.38109484726567
$ws.Range("I23").Value = 19.78462438174012
$ws.Range("J23").Value = 9.136767751826069
$ws.Range("M23").Value = 17.98649228774777
$ws.Range("N23").Value = 17.18608209020867
$ws.Range("O23").Value = 22.18767798997643
$ws.Range("B24").Value = 13.98581382768327
$ws.Range("C24").Value = 8.930579128612976
$ws.Range("D24").Value = 14.62622447434662
$ws.Range("E24").Value = 15.7262439191285
$ws.Range("G24").Value = 30.27972633340427
$ws.Range("H24").Value = 14.42612966897375
$ws.Range("I24").Value = 19.90876220850854
$ws.Range("J24").Value = 9.167167832516601
$ws.Range("M24").Value = 17.70104278465137
$ws.Range("N24").Value = 17.25312866436225
$ws.Range("O24").Value = 22.21476765910668
$ws.Range("B25").Value = 13.10507293717164
$ws.Range("C25").Value = 8.205119687651033
$ws.Range("D25").Value = 14.56966102535404
$ws.Range("E25").Value = 15.73661428093481
$ws.Range("G25").Value = 30.13433675788395
$ws.Range("H25").Value = 14.48461883037671
$ws.Range("I25").Value = 20.05892930041602
$ws.Range("J25").Value = 9.203283490836533
$ws.Range("M25").Value = 17.39738132809608
$ws.Range("N25").Value = 17.33187394586058
$ws.Range("O25").Value = 22.2644684393515
